$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.189.84"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.597.16"
$ws.Range("E3").Value = "  +3.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.990.95"
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.108"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.518.44"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "43.304.87"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.92%  "
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0809"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.81%  "
$ws.Range("E36").Value = "  +3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("E39").Value = "  +9.78%  "
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("E42").Value = "  +6.02%  "
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").Value = "2.019.43"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "83.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("E49").Value = "  +4.96%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.789.17"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.61%  "
